# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Re-order Ecuador / Chile / Luxemburgo (rows 29-31) and refresh Ecuador + Japon data ---
# Row 29 becomes Ecuador with updated figures
$ws.Cells.Item(29, 1).Value = "Ecuador"
$ws.Cells.Item(29, 2).Value = 1627
$ws.Cells.Item(29, 3).Value = 224
$ws.Cells.Item(29, 4).Value = 3
$ws.Cells.Item(29, 5).Value = 1583
$ws.Cells.Item(29, 6).Value = 58
$ws.Cells.Item(29, 7).Value = 7
$ws.Cells.Item(29, 8).Value = 41

# Row 30 becomes Chile, keeping its previous figures
$ws.Cells.Item(30, 1).Value = "Chile"
$ws.Cells.Item(30, 2).Value = 1610
$ws.Cells.Item(30, 3).Value = 304
$ws.Cells.Item(30, 4).Value = 43
$ws.Cells.Item(30, 5).Value = 1562
$ws.Cells.Item(30, 6).Value = 7
$ws.Cells.Item(30, 7).Value = 1
$ws.Cells.Item(30, 8).Value = 5

# Row 31 becomes Luxemburgo, keeping its previous figures
$ws.Cells.Item(31, 1).Value = "Luxemburgo"
$ws.Cells.Item(31, 2).Value = 1605
$ws.Cells.Item(31, 3).Value = 152
$ws.Cells.Item(31, 4).Value = 40
$ws.Cells.Item(31, 5).Value = 1550
$ws.Cells.Item(31, 6).Value = 25
$ws.Cells.Item(31, 7).Value = 6
$ws.Cells.Item(31, 8).Value = 15

# Row 32 stays Japon, figures refreshed
$ws.Cells.Item(32, 2).Value = 1499
$ws.Cells.Item(32, 3).Value = 112
$ws.Cells.Item(32, 5).Value = 1078

# --- Re-order Venezuela / Afganistan (rows 99-100) ---
# Row 99 becomes Venezuela with updated figures
$ws.Cells.Item(99, 1).Value = "Venezuela"
$ws.Cells.Item(99, 2).Value = 113
$ws.Cells.Item(99, 3).Value = 6
$ws.Cells.Item(99, 4).Value = 31
$ws.Cells.Item(99, 5).Value = 80
$ws.Cells.Item(99, 6).Value = 2
$ws.Cells.Item(99, 7).Value = 1
$ws.Cells.Item(99, 8).Value = 2

# Row 100 becomes Afganistan, keeping its previous figures
$ws.Cells.Item(100, 1).Value = "Afganistan"
$ws.Cells.Item(100, 2).Value = 110
$ws.Cells.Item(100, 3).Value = 16
$ws.Cells.Item(100, 4).Value = 2
$ws.Cells.Item(100, 5).Value = 104
$ws.Cells.Item(100, 6).Value = 0
$ws.Cells.Item(100, 7).Value = 0
$ws.Cells.Item(100, 8).Value = 4

# --- Plain figure refreshes (country stays the same) ---
# Estados Unidos
$ws.Cells.Item(4, 2).Value = 102464
$ws.Cells.Item(4, 3).Value = 17029
$ws.Cells.Item(4, 4).Value = 2471
$ws.Cells.Item(4, 5).Value = 98386
$ws.Cells.Item(4, 7).Value = 312
$ws.Cells.Item(4, 8).Value = 1607

# China
$ws.Cells.Item(6, 2).Value = 81394
$ws.Cells.Item(6, 3).Value = 54
$ws.Cells.Item(6, 4).Value = 74971
$ws.Cells.Item(6, 5).Value = 3128
$ws.Cells.Item(6, 6).Value = 886
$ws.Cells.Item(6, 7).Value = 3
$ws.Cells.Item(6, 8).Value = 3295

# Alemania
$ws.Cells.Item(8, 5).Value = 43862
$ws.Cells.Item(8, 6).Value = 1581
$ws.Cells.Item(8, 7).Value = 84
$ws.Cells.Item(8, 8).Value = 351

# Austria
$ws.Cells.Item(15, 2).Value = 7697
$ws.Cells.Item(15, 3).Value = 788
$ws.Cells.Item(15, 5).Value = 7414

# Canada
$ws.Cells.Item(18, 2).Value = 4757
$ws.Cells.Item(18, 3).Value = 714
$ws.Cells.Item(18, 4).Value = 354
$ws.Cells.Item(18, 5).Value = 4348
$ws.Cells.Item(18, 7).Value = 16
$ws.Cells.Item(18, 8).Value = 55

# Noruega
$ws.Cells.Item(20, 2).Value = 3771
$ws.Cells.Item(20, 3).Value = 399
$ws.Cells.Item(20, 4).Value = 7
$ws.Cells.Item(20, 5).Value = 3745

# Irlanda
$ws.Cells.Item(27, 6).Value = 59

# Tunez
$ws.Cells.Item(81, 5).Value = 218
$ws.Cells.Item(81, 7).Value = 2
$ws.Cells.Item(81, 8).Value = 7

# --- Timestamp footer ---
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 28 de Marzo de 2020 a las 01:29"
